$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "66.867.60"
$ws.Range("E2").Value = "  +1.28%  "
$ws.Range("D3").Value = "3.904.19"
$ws.Range("E3").Value = "  +3.00%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "'467.41"
$ws.Range("E5").Value = "  +9.54%  "
$ws.Range("D6").Value = "'145.28"
$ws.Range("E6").Value = "  +6.14%  "
$ws.Range("E7").Value = "  +2.67%  "
$ws.Range("D8").Value = "'0.999"
$ws.Range("E8").Value = "  +0.00%  "
$ws.Range("D9").Value = "'0.740"
$ws.Range("E9").Value = "  -0.17%  "
$ws.Range("D10").Value = "'0.164"
$ws.Range("E10").Value = "  +3.35%  "
$ws.Range("D11").Value = "'0.0000337"
$ws.Range("E11").Value = "  +2.10%  "
$ws.Range("D12").Value = "'43.25"
$ws.Range("E12").Value = "  +0.00%  "
$ws.Range("D13").Value = "'10.45"
$ws.Range("E13").Value = "  -2.21%  "
$ws.Range("D14").Value = "4.527.82"
$ws.Range("E14").Value = "  +3.45%  "
$ws.Range("D15").Value = "'15.03"
$ws.Range("E15").Value = "  -0.30%  "
$ws.Range("D16").Value = "3.883.25"
$ws.Range("E16").Value = "  +1.99%  "
$ws.Range("E17").Value = "  -0.38%  "
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("E19").Value = "  +2.89%  "
$ws.Range("D20").Value = "67.172.27"
$ws.Range("E20").Value = "  +1.68%  "
$ws.Range("D21").Value = "'431.51"
$ws.Range("E21").Value = "  +4.69%  "
$ws.Range("D22").Value = "'14.71"
$ws.Range("E22").Value = "  -3.53%  "
$ws.Range("D23").Value = "'3.35"
$ws.Range("E23").Value = "  +1.78%  "
$ws.Range("E24").Value = "  +3.21%  "
$ws.Range("D25").Value = "'38.67"
$ws.Range("E25").Value = "  +4.23%  "
$ws.Range("D26").Value = "'3.52"
$ws.Range("E26").Value = "  +6.38%  "
$ws.Range("B27").Value = "LEO"
$ws.Range("C27").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D27").Value = "'5.67"
$ws.Range("E27").Value = "  +5.07%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").Value = "'10.10"
$ws.Range("E28").Value = "  +2.72%  "
$ws.Range("D29").Value = "'9.73"
$ws.Range("E29").Value = "  -2.24%  "
$ws.Range("D30").Value = "'737.31"
$ws.Range("E30").Value = "  +4.54%  "
$ws.Range("D31").Value = "'13.69"
$ws.Range("E31").Value = "  -2.87%  "
$ws.Range("E32").Value = "  +2.13%  "
$ws.Range("E33").Value = "  -0.07%  "
$ws.Range("D34").Value = "'43.88"
$ws.Range("E34").Value = "  +9.85%  "
$ws.Range("E35").Value = "  +3.19%  "
$ws.Range("E36").Value = "  +3.71%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("D38").Value = "'3.26"
$ws.Range("E38").Value = "  +12.72%  "
$ws.Range("D39").Value = "'5.37"
$ws.Range("E39").Value = "  -7.68%  "
$ws.Range("D40").Value = "'0.0477"
$ws.Range("E40").Value = "  +0.49%  "
$ws.Range("D41").Value = "0.0₃0754"
$ws.Range("E41").Value = "  +9.19%  "
$ws.Range("E42").Value = "  +0.90%  "
$ws.Range("E43").Value = "  +4.19%  "
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  +6.23%  "
$ws.Range("E46").Value = "  +5.58%  "
$ws.Range("E47").Value = "  +1.50%  "
$ws.Range("D48").Value = "'2.49"
$ws.Range("E48").Value = "  -5.98%  "
$ws.Range("D49").Value = "'3.18"
$ws.Range("E49").Value = "  +0.46%  "
$ws.Range("B50").Value = "Stacks"
$ws.Range("C50").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D50").Value = "'2.88"
$ws.Range("E50").Value = "  +1.92%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'26.41"
$ws.Range("E51").Value = "  +2.90%  "
